# Apply "everything works but lowercase class names on del" change:
# admin10 (row 5 on "Timeslot Information") gets enrolled in English (Monday),
# Java (Tuesday) and Python (Thursday); the corresponding "Current Occupancy"
# counters on "Class Information" are bumped accordingly.

$wb = $excel.ActiveWorkbook

$wsTimeslot = $wb.Worksheets.Item("Timeslot Information")
$wsTimeslot.Range("B5").Value = "English,null,null,null,Biology,null,null,null,null,null,Java,null"
$wsTimeslot.Range("C5").Value = "null,null,null,null,null,null,null,null,null,null,Java,null"
$wsTimeslot.Range("E5").Value = "null,null,null,null,null,null,null,null,null,Python,null,null"

$wsClass = $wb.Worksheets.Item("Class Information")
$wsClass.Cells.Item(2, 7).Value = 1   # English  -> Current Occupancy
$wsClass.Cells.Item(11, 7).Value = 2  # Python   -> Current Occupancy
$wsClass.Cells.Item(12, 7).Value = 2  # Java     -> Current Occupancy
